# Updated capital structure database
# Refresh the computed ratio/metric columns (D, G:R, U:AQ) for the
# Lithuania engineering/construction data rows (rows 2 and 3) to reflect
# the latest source figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = -0.0545
$ws.Range("G2").Value = 0.006754068716094033
$ws.Range("H2").Value = 0.006754068716094033
$ws.Range("I2").Value = -0.02305605786618445
$ws.Range("J2").Value = -0.02305605786618445
$ws.Range("K2").Value = -13.9
$ws.Range("L2").Value = -0.1256781193490054
$ws.Range("M2").Value = 0.549
$ws.Range("N2").Value = 0.04815789473684211
$ws.Range("O2").Value = -0.03949640287769784
$ws.Range("P2").Value = 0.549
$ws.Range("Q2").Value = 0.04815789473684211
$ws.Range("R2").Value = -0.03949640287769784
$ws.Range("U2").Value = 5.35
$ws.Range("V2").Value = 0.469298245614035
$ws.Range("W2").Value = -0.3518987341772152
$ws.Range("X2").Value = 0.1132496740213998
$ws.Range("Y2").Value = -0.465148408198615
$ws.Range("Z2").Value = 2.642140468227425
$ws.Range("AA2").Value = -0.06091734352603918
$ws.Range("AB2").Value = 0.06001289839764713
$ws.Range("AC2").Value = -0.1209302419236863
$ws.Range("AD2").Value = 17.5
$ws.Range("AF2").Value = 17.5
$ws.Range("AG2").Value = 12.15
$ws.Range("AH2").Value = 0.6055363321799309
$ws.Range("AI2").Value = 0.3995433789954338
$ws.Range("AJ2").Value = 0.5159235668789809
$ws.Range("AK2").Value = 0.3159947984395318
$ws.Range("AL2").Value = 4.98
$ws.Range("AM2").Value = 4.98
$ws.Range("AN2").Value = -18.79699248120301
$ws.Range("AO2").Value = -0.5120481927710843
$ws.Range("AP2").Value = -13.05048335123523
$ws.Range("AQ2").Value = -0.5120481927710843

# Row 3
$ws.Range("D3").Value = -0.0545
$ws.Range("G3").Value = 0.006754068716094033
$ws.Range("H3").Value = 0.006754068716094033
$ws.Range("I3").Value = -0.02305605786618445
$ws.Range("J3").Value = -0.02305605786618445
$ws.Range("K3").Value = -13.9
$ws.Range("L3").Value = -0.1256781193490054
$ws.Range("M3").Value = 0.549
$ws.Range("N3").Value = 0.04815789473684211
$ws.Range("O3").Value = -0.03949640287769784
$ws.Range("P3").Value = 0.549
$ws.Range("Q3").Value = 0.04815789473684211
$ws.Range("R3").Value = -0.03949640287769784
$ws.Range("U3").Value = 5.35
$ws.Range("V3").Value = 0.469298245614035
$ws.Range("W3").Value = -0.3518987341772152
$ws.Range("X3").Value = 0.1132496740213998
$ws.Range("Y3").Value = -0.465148408198615
$ws.Range("Z3").Value = 2.642140468227425
$ws.Range("AA3").Value = -0.06091734352603918
$ws.Range("AB3").Value = 0.06001289839764713
$ws.Range("AC3").Value = -0.1209302419236863
$ws.Range("AD3").Value = 17.5
$ws.Range("AF3").Value = 17.5
$ws.Range("AG3").Value = 12.15
$ws.Range("AH3").Value = 0.6055363321799309
$ws.Range("AI3").Value = 0.3995433789954338
$ws.Range("AJ3").Value = 0.5159235668789809
$ws.Range("AK3").Value = 0.3159947984395318
$ws.Range("AL3").Value = 4.98
$ws.Range("AM3").Value = 4.98
$ws.Range("AN3").Value = -18.79699248120301
$ws.Range("AO3").Value = -0.5120481927710843
$ws.Range("AP3").Value = -13.05048335123523
$ws.Range("AQ3").Value = -0.5120481927710843
